$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to hold a literal text string even when $val looks
    # like a plain decimal number (e.g. "212.68"), matching the workbook's
    # inlineStr/text storage for the Price column. NumberFormat "@" makes
    # Excel store the assignment as text instead of coercing it to a
    # number; resetting the style back to Normal afterwards avoids leaving
    # a lasting number-format override on the cell.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.548.51"
$ws.Range("E2").Value = "  +0.51%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.641.30"
$ws.Range("E3").Value = "  -0.85%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5 - BNB
Set-TextValue "D5" "212.68"
$ws.Range("E5").Value = "  -0.37%  "

# Row 6 - XRP
Set-TextValue "D6" "0.535"
$ws.Range("E6").Value = "  +4.37%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.02%  "

# Row 8 - Solana
Set-TextValue "D8" "23.01"
$ws.Range("E8").Value = "  -4.15%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -1.70%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.60%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +1.36%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.873.80"
$ws.Range("E12").Value = "  -0.80%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.623.34"
$ws.Range("E13").Value = "  -1.90%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -1.05%  "

# Row 15 - Polygon
Set-TextValue "D15" "0.563"
$ws.Range("E15").Value = "  -2.14%  "

# Row 16 - Litecoin
Set-TextValue "D16" "64.03"
$ws.Range("E16").Value = "  -2.64%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "27.545.64"
$ws.Range("E17").Value = "  +0.47%  "

# Row 18 - BitcoinCash
Set-TextValue "D18" "229.12"
$ws.Range("E18").Value = "  -1.29%  "

# Row 19 - ShibaInu
$ws.Range("E19").Value = "  -0.35%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  +1.55%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.00%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -1.76%  "

# Row 23 - Avalanche
Set-TextValue "D23" "10.00"
$ws.Range("E23").Value = "  +7.32%  "

# Row 24 - Toncoin
Set-TextValue "D24" "1.94"
$ws.Range("E24").Value = "  -3.75%  "

# Row 25 - Monero
Set-TextValue "D25" "149.49"
$ws.Range("E25").Value = "  +1.72%  "

# Row 26 - Cosmos
Set-TextValue "D26" "6.97"
$ws.Range("E26").Value = "  -3.44%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  +1.36%  "

# Row 28 - BinanceUSD
$ws.Range("E28").Value = "  -0.12%  "

# Row 29 - EthereumClassic
Set-TextValue "D29" "15.60"
$ws.Range("E29").Value = "  -2.01%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -1.01%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  -2.23%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -0.02%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextValue "D33" "3.17"
$ws.Range("E33").Value = "  +1.50%  "

# Row 34 - Maker
$ws.Range("D34").Value = "1.428.04"
$ws.Range("E34").Value = "  -2.63%  "

# Row 35 - LidoDAOToken
Set-TextValue "D35" "1.59"

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  -1.97%  "

# Row 37 - ImmutableX
Set-TextValue "D37" "0.574"
$ws.Range("E37").Value = "  +0.20%  "

# Row 38 - ARBITRUM
Set-TextValue "D38" "0.876"
$ws.Range("E38").Value = "  -3.99%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  -1.44%  "

# Row 40 - TrustWalletToken
Set-TextValue "D40" "0.900"
$ws.Range("E40").Value = "  +15.08%  "

# Row 41 - WEMIXToken
$ws.Range("E41").Value = "  -1.71%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  +0.03%  "

# Row 43 - mCoin
$ws.Range("E43").Value = "  -1.17%  "

# Row 44 - FraxShare
$ws.Range("E44").Value = "  +1.21%  "

# Row 45 - MXToken
$ws.Range("E45").Value = "  +1.62%  "

# Row 46 - Aave
Set-TextValue "D46" "65.19"
$ws.Range("E46").Value = "  -0.21%  "

# Row 47 - RocketPoolETH
$ws.Range("D47").Value = "1.783.13"
$ws.Range("E47").Value = "  -0.75%  "

# Row 48 - RenderToken
$ws.Range("E48").Value = "  -2.29%  "

# Row 49 - Quant
Set-TextValue "D49" "86.32"
$ws.Range("E49").Value = "  -2.28%  "

# Row 50 - BabyDogeCoin
$ws.Range("E50").Value = "  +0.94%  "

# Row 51 - Algorand
$ws.Range("E51").Value = "  -2.71%  "
